$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sampletext"
$ws.Range("B1").Value = "hello"
$ws.Range("B2").Value = "There"
$ws.Range("B3").Value = "illegal"
$ws.Range("B4").Value = "will fail"
$ws.Range("B5").Value = "text"
$ws.Range("B6").Value = "generate"
$ws.Range("B7").Value = "document"
$ws.Range("B8").Value = "legal"

$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
